# Fix project gantt chart
#
# - Insert a "ProjectID" column right after "ID" (new column B), giving each
#   task row the id of the project it belongs to.
# - Rename "Dependencies" -> "TaskDependencies".
# - Insert a new "ProjectDependency" column right after "TaskDependencies",
#   flagging tasks whose project depends on another project.
# - ProjectName/TaskName/EstimatedEffortHours/Progress all shift right to make
#   room for the two new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Physically shift the existing columns out of the way first (this keeps the
# untouched columns' widths/number-formats intact) and THEN grow the table to
# cover the new A1:H6 range.
$ws.Columns("B").Insert()   # make room for ProjectID right after ID
$ws.Columns("G").Insert()   # make room for ProjectDependency after TaskDependencies
$lo.Resize($ws.Range("A1:H6"))

# --- Header row (row 1) : final column order ---
# A=ID  B=ProjectID  C=ProjectName  D=TaskName  E=EstimatedEffortHours
# F=TaskDependencies  G=ProjectDependency  H=Progress
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "ProjectID"
$ws.Range("C1").Value = "ProjectName"
$ws.Range("D1").Value = "TaskName"
$ws.Range("E1").Value = "EstimatedEffortHours"
$ws.Range("F1").Value = "TaskDependencies"
$ws.Range("G1").Value = "ProjectDependency"
$ws.Range("H1").Value = "Progress"

# --- Data rows ---
# Row 2: ID=1, ProjectID=1, ProjectName=Project A, TaskName=Design UI, Hours=100, Progress=50
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Project A"
$ws.Range("D2").Value = "Design UI"
$ws.Range("E2").Value = 100
$ws.Range("H2").Value = 50

# Row 3: ID=2, ProjectID=1, ProjectName=Project A, TaskName=Implement Backend, Hours=55, TaskDependencies=1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Project A"
$ws.Range("D3").Value = "Implement Backend"
$ws.Range("E3").Value = 55
$ws.Range("F3").Value = 1

# Row 4: ID=3, ProjectID=1, ProjectName=Project A, TaskName=Testing, Hours=15, TaskDependencies=2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Project A"
$ws.Range("D4").Value = "Testing"
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 2

# Row 5: ID=4, ProjectID=2, ProjectName=Project B, TaskName=Database Setup, Hours=60, ProjectDependency=1
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Project B"
$ws.Range("D5").Value = "Database Setup"
$ws.Range("E5").Value = 60
$ws.Range("G5").Value = 1

# Row 6: ID=5, ProjectID=2, ProjectName=Project B, TaskName=API Development, Hours=50, TaskDependencies=4, ProjectDependency=1
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Project B"
$ws.Range("D6").Value = "API Development"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1

# --- Column widths for the two brand-new columns (B, G) ---
# (A, C, D, E, F, H keep the exact widths/best-fit flags they already had,
# since they were shifted into place rather than rewritten.)
$ws.Columns("B").ColumnWidth = 12
$ws.Columns("G").ColumnWidth = 14.5

# --- Selection ---
$ws.Range("G7").Select()
